$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value to a cell as literal TEXT (not auto-converted to a
# number by the usual `.Value =` coercion), while keeping the cell's style
# unchanged. We briefly force a "Text" number format so the typed value is
# stored verbatim, then paste back just the *formats* from a same-sheet donor
# cell that already carries the style we want to end up with.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($cell, $donor, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $donor.Copy()
    $cell.PasteSpecial(-4122)  # xlPasteFormats
}

$totalSheet = $wb.Worksheets.Item("总计")
$q1Sheet = $wb.Worksheets.Item("2022-Q1")
$q32021Sheet = $wb.Worksheets.Item("2021-Q3")

# ---------------------------------------------------------------------------
# Build the new "2022-Q3" sheet by duplicating "2022-Q1" (same column
# layout/styling) right before it, then overwrite its contents in place.
# This keeps every inherited style (bold+bordered header, index-column
# style, etc.) byte-identical to the sibling quarter sheets.
# ---------------------------------------------------------------------------
$q1Sheet.Copy($q1Sheet)
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# "2022-Q1" has 3 data rows; "2022-Q3" only needs 2, so drop the extra row.
$q3Sheet.Rows(4).Delete()

$q3TextDonor = $q3Sheet.Cells.Item(1, 1)  # untouched, default-style (0) cell used purely as a format donor

Set-TextValue $q3Sheet.Cells.Item(2, 2) $q3TextDonor "007368"
Set-TextValue $q3Sheet.Cells.Item(2, 3) $q3TextDonor "浙商沪港深精选混合A"
Set-TextValue $q3Sheet.Cells.Item(2, 4) $q3TextDonor "6.59"
Set-TextValue $q3Sheet.Cells.Item(2, 5) $q3TextDonor "84.00"
Set-TextValue $q3Sheet.Cells.Item(2, 6) $q3TextDonor "4.80"
Set-TextValue $q3Sheet.Cells.Item(2, 7) $q3TextDonor "0.3163"
$q3Sheet.Cells.Item(2, 8).Value = 9

Set-TextValue $q3Sheet.Cells.Item(3, 2) $q3TextDonor "007369"
Set-TextValue $q3Sheet.Cells.Item(3, 3) $q3TextDonor "浙商沪港深精选混合C"
Set-TextValue $q3Sheet.Cells.Item(3, 4) $q3TextDonor "0.32"
Set-TextValue $q3Sheet.Cells.Item(3, 5) $q3TextDonor "84.00"
Set-TextValue $q3Sheet.Cells.Item(3, 6) $q3TextDonor "4.80"
Set-TextValue $q3Sheet.Cells.Item(3, 7) $q3TextDonor "0.0154"
$q3Sheet.Cells.Item(3, 8).Value = 9

# ---------------------------------------------------------------------------
# "总计" sheet: insert a new row 2 for "2022-Q3", shifting the existing
# "2022-Q1" row down to row 3 and "2021-Q3" row down to row 4.
# ---------------------------------------------------------------------------
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(4, 2).Value = "2021-Q3"
$totalSheet.Cells.Item(4, 3).Value = 1
$totalSheet.Cells.Item(4, 4).Value = 0.04
$totalSheet.Cells.Item(3, 1).Copy()
$totalSheet.Cells.Item(4, 1).PasteSpecial(-4122)  # xlPasteFormats
$totalSheet.Cells.Item(4, 1).Value = 2

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(3, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(3, 3).Value = 3
$totalSheet.Cells.Item(3, 4).Value = 0.14

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 2
$totalSheet.Cells.Item(2, 4).Value = 0.33

# Restore "2021-Q3" as the active sheet/tab (unchanged from the original file).
# (Re-look the sheet up by name rather than reusing the earlier variable: sheet
# references here are positional, and inserting/renaming sheets above it would
# otherwise leave the old variable pointing at the wrong sheet.)
$wb.Worksheets.Item("2021-Q3").Activate()

Write-Output "done"
